$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the two new columns
$ws.Range("F1").Value = "Battle Image"
$ws.Range("G1").Value = "Sprite Image"

# Per-row sprite/battle image file names (first few rows have distinct art,
# the remainder currently reuse the placeholder 000.png for both columns).
# Rows are populated in the same order they were authored so new shared
# strings land at the same indices as the source edit.
$ws.Range("F2").Value = "Engimon/004.png"
$ws.Range("G2").Value = "Engimon/000.png"

$ws.Range("F6").Value = "Engimon/007.png"
$ws.Range("G6").Value = "Engimon/000.png"

$ws.Range("F4").Value = "Engimon/011.png"
$ws.Range("G4").Value = "Engimon/000.png"

$ws.Range("F5").Value = "Engimon/053.png"
$ws.Range("G5").Value = "Engimon/000.png"

$ws.Range("F3").Value = "Engimon/333.png"
$ws.Range("G3").Value = "Engimon/000.png"

for ($r = 7; $r -le 37; $r++) {
    $ws.Range("F$r").Value = "Engimon/000.png"
    $ws.Range("G$r").Value = "Engimon/000.png"
}

# Column widths for the newly-used / resized columns (closest values the
# ColumnWidth property's internal rounding can land on the target widths)
$ws.Columns.Item(5).ColumnWidth = 31.33
$ws.Columns.Item(6).ColumnWidth = 20.66
$ws.Columns.Item(7).ColumnWidth = 20.66

# Match the saved selection from the edit
$ws.Range("F6").Select() | Out-Null
